$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 231 (pushes existing rows 231:265 down to 232:266,
# carrying their formatting/styles with them).
$ws.Rows.Item(231).Insert()

# Populate the new row 231 with the same record as the old row 231
# (A..I, K..R unchanged) except Fecha (D) and Volumen (J), which take
# the new values.
$ws.Range("A231").Value = 10
$ws.Range("B231").Value = "Vega Modelo de Temuco"
$ws.Range("C231").Value = "La Araucanía"
$ws.Range("D231").Value = 44474
$ws.Range("E231").Value = 9
$ws.Range("F231").Value = 100112023
$ws.Range("G231").Value = "Brócoli"
$ws.Range("H231").Value = "Sin especificar"
$ws.Range("I231").Value = "Primera"
$ws.Range("J231").Value = 1500
$ws.Range("K231").Value = 800
$ws.Range("L231").Value = 800
$ws.Range("M231").Value = 800
$ws.Range("N231").Value = "$/unidad"
$ws.Range("O231").Value = "Región Metropolitana"
$ws.Range("P231").Value = 800
$ws.Range("Q231").Value = 1
$ws.Range("R231").Value = "Hortaliza"
